$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1, matching the style of the existing H1 header
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF) for rows 2-79
$iVals = @(9,8,8,9,8,9,9,8,6,9,6,8,8,9,9,9,9,9,7,9,9,9,9,10,9,9,9,9,8,9,8,8,9,10,9,9,8,9,8,9,9,9,9,8,8,9,9,10,9,9,9,8,9,8,8,8,8,9,9,9,9,6,9,9,9,9,8,9,7,5,7,7,9,7,5,6,5,4)
$jVals = @(9,9,8,9,8,9,10,8,7,9,6,8,8,9,9,9,9,9,7,9,9,9,9,10,9,9,9,9,9,9,9,8,9,10,9,9,9,9,8,9,9,9,9,9,8,9,9,10,9,9,9,9,9,8,9,8,8,9,9,9,10,6,9,9,9,9,8,9,7,7,7,7,9,7,6,6,5,4)

$startRow = 2
for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
  $r = $startRow + $idx
  $ws.Cells.Item($r, 9).Value = $iVals[$idx]
  $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
